# Sept 2 afternoon update
# Applies the Sept 2 afternoon data-entry update to the "Kelp consumption"
# sheet (urchin size / weight-before / weight-after readings for day 2,
# plus trial end times and a few "Percent of Kelp Consumed" placeholders),
# fixes one weight reading, adds a couple of literal "Times Crossing Cory"
# counts + a comment on the "Videos" sheet, and adjusts window/selection
# state.

$wb = $excel.ActiveWorkbook

$kelp = $wb.Worksheets.Item("Kelp consumption")
$videos = $wb.Worksheets.Item("Videos")

# --- one-off correction -----------------------------------------------
$kelp.Range("M13").Value = 43

# --- "Percent of Kelp Consumed" placeholders (col T) added as 0 -------
foreach ($r in 3,8,9,10,11,15,16,18,19,21,23,24,25,27,28,29) {
    $kelp.Range("T" + $r).Value = 0
}

# --- newly recorded urchin size (K) / weight before (L) / weight
#     after (M) / stop-trial clock (O) / stop time (P) readings --------
$rows = @{
    16 = @{ K = 45; L = 45; M = 45; O = 1910; P = 710 }
    17 = @{ K = 43; L = 50; M = 41; O = 1910; P = 710 }
    18 = @{ K = 58; L = 85; M = 85; O = 1910; P = 710 }
    19 = @{ K = 48; L = 51; M = 51; O = 1910; P = 710 }
    20 = @{ K = 50; L = 60; M = 61; O = 1910; P = 710 }
    21 = @{ K = 49; L = 46; M = 45; O = 1910; P = 710 }
    22 = @{ K = 48; L = 54; M = 53; O = 1905; P = 705 }
    23 = @{ K = 57; L = 79; M = 79; O = 1905; P = 705 }
    24 = @{ K = 50; L = 50; M = 51; O = 1905; P = 705 }
    25 = @{ K = 50; L = 55; M = 56; O = 1905; P = 705 }
    26 = @{ K = 61; L = 92; M = 94; O = 1905; P = 705 }
    27 = @{ K = 44; L = 34; M = 35; O = 1905; P = 705 }
    28 = @{ K = 45; L = 44; M = 45; O = 1905; P = 705 }
    29 = @{ K = 47; L = 40; M = 41; O = 1905; P = 705 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $kelp.Range("K$r").Value = $vals.K
    $kelp.Range("L$r").Value = $vals.L
    $kelp.Range("M$r").Value = $vals.M
    $kelp.Range("O$r").Value = $vals.O
    $kelp.Range("P$r").Value = $vals.P
}

# Row 16's K/L/M cells were previously blank-but-bold placeholders;
# drop the bold now that they hold real data.
$kelp.Range("K16:M16").Font.Bold = $false

# --- "Videos" sheet: literal "Times Crossing Cory" counts (col Q) -----
$videos.Range("Q2").Value = 5
$videos.Range("Q3").Value = 0
$videos.Range("Q4").Value = 1
$videos.Range("Q5").Value = 1
$videos.Range("Q6").Value = 5
$videos.Range("Q7").Value = 8

# New comment for the trial whose GoPro lost its BacPack connection
$videos.Range("X4").Value = "GoPro not connected to BacPack"

# --- view / selection state --------------------------------------------
$kelp.Activate()
$kelp.Application.ActiveWindow.ScrollColumn = 8   # topLeftCell="H1"
$kelp.Range("N24").Select()

$videos.Activate()
$videos.Application.ActiveWindow.ScrollRow = 1    # topLeftCell back to A1
$videos.Range("R3").Select()

$kelp.Activate()
